$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Workbook-level defined name left behind by the "MySQL for Excel" add-in
# (picked up automatically the next time the file was opened/saved).
$mysqlDateFormat = $wb.Names.Add("LOCAL_MYSQL_DATE_FORMAT", "=REPT(LOCAL_YEAR_FORMAT,4)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_MONTH_FORMAT,2)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_DAY_FORMAT,2)&"" ""&REPT(LOCAL_HOUR_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_MINUTE_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_SECOND_FORMAT,2)")
$mysqlDateFormat.Visible = $false

# Append a new "Icon" property row (row 19) to the Property table.
# Shared-string intern order matters: "显示图标" must be registered before
# "Icon" to match the target uniqueCount ordering, so set J before A.
$row = 19
$ws.Range("J$row").Value = "显示图标"
$ws.Range("A$row").Value = "Icon"
$ws.Range("B$row").Value = "string"
$ws.Range("C$row").Value = $false
$ws.Range("D$row").Value = $false
$ws.Range("E$row").Value = $false
$ws.Range("F$row").Value = $true
$ws.Range("G$row").Value = 0
$ws.Range("H$row").Value = 0
$ws.Range("I$row").Value = "Friend"

# Match the text format ("@") used by columns A, B, I, J elsewhere in the table.
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("B$row").NumberFormat = "@"
$ws.Range("I$row").NumberFormat = "@"
$ws.Range("J$row").NumberFormat = "@"

# The TRUE/FALSE list validation previously had a gap at row 19
# (F2:F18 + F20:F1048576); now that row 19 is filled in, Excel re-merges
# it into one contiguous range.
$ws.Range("F2:F1048576").Validation.Delete()
$ws.Range("F2:F1048576").Validation.Add(3, 1, 1, '"TRUE,FALSE"')

[void]$ws.Range("H24").Select()
